# Vincula un profesor con un gasto promocional:
# actualiza las horas consumidas del "Dia 1" (columna T) de la tarea en
# la fila 13 de la hoja "Casos de Uso", de 1 a 2 horas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Actualiza la celda de horas consumidas (Dia 1) de 1 a 2
$ws.Range("T13").Value = 2

# Deja seleccionada la columna N completa, como quedo al guardar el archivo
$ws.Range("N1:N1048576").Select() | Out-Null
